$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.1079986666666667
$ws.Range("H2").Value = 0.323996
$ws.Range("I2").Value = 0.004187739561209694
$ws.Range("J2").Value = 0.004187739561209694
$ws.Range("M2").Value = 10.34761366666667
$ws.Range("N2").Value = 31.042841
$ws.Range("O2").Value = 0.2299953477621856
$ws.Range("P2").Value = 0.2299953477621856
$ws.Range("Q2").Value = 1.117528479181778
$ws.Range("R2").Value = 10.057756312636
$ws.Range("S2").Value = 0.000963160616717886
$ws.Range("T2").Value = 0.0009631606167178862

$ws.Range("G3").Value = 0.1079986666666667
$ws.Range("H3").Value = 0.323996
$ws.Range("I3").Value = 0.004187739561209694
$ws.Range("J3").Value = 0.004187739561209694
$ws.Range("O3").Value = 0.6794731949692173
$ws.Range("P3").Value = 0.6794731949692174
$ws.Range("Q3").Value = 3.301504372183556
$ws.Range("R3").Value = 29.713539349652
$ws.Range("S3").Value = 0.002845456779354139
$ws.Range("T3").Value = 0.00284545677935414

$ws.Range("G4").Value = 0.1079986666666667
$ws.Range("H4").Value = 0.323996
$ws.Range("I4").Value = 0.004187739561209694
$ws.Range("J4").Value = 0.004187739561209694
$ws.Range("M4").Value = 4.073058666666666
$ws.Range("N4").Value = 12.219176
$ws.Range("O4").Value = 0.09053145726859702
$ws.Range("P4").Value = 0.09053145726859703
$ws.Range("Q4").Value = 0.4398849052551111
$ws.Range("R4").Value = 3.958964147296
$ws.Range("S4").Value = 0.0003791221651376686
$ws.Range("T4").Value = 0.0003791221651376687

$ws.Range("I5").Value = 0.9687110856121154
$ws.Range("J5").Value = 0.9687110856121155
$ws.Range("M5").Value = 10.34761366666667
$ws.Range("N5").Value = 31.042841
$ws.Range("O5").Value = 0.2299953477621856
$ws.Range("P5").Value = 0.2299953477621856
$ws.Range("Q5").Value = 258.5075338252222
$ws.Range("R5").Value = 2326.567804427
$ws.Range("S5").Value = 0.2227990430164428
$ws.Range("T5").Value = 0.2227990430164428

$ws.Range("I6").Value = 0.9687110856121154
$ws.Range("J6").Value = 0.9687110856121155
$ws.Range("O6").Value = 0.6794731949692173
$ws.Range("P6").Value = 0.6794731949692174
$ws.Range("S6").Value = 0.658213216342963
$ws.Range("T6").Value = 0.6582132163429633

$ws.Range("I7").Value = 0.9687110856121154
$ws.Range("J7").Value = 0.9687110856121155
$ws.Range("M7").Value = 4.073058666666666
$ws.Range("N7").Value = 12.219176
$ws.Range("O7").Value = 0.09053145726859702
$ws.Range("P7").Value = 0.09053145726859703
$ws.Range("Q7").Value = 101.7545092968889
$ws.Range("R7").Value = 915.790583672
$ws.Range("S7").Value = 0.08769882625270944
$ws.Range("T7").Value = 0.08769882625270947

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.6989190000000001
$ws.Range("H8").Value = 2.096757
$ws.Range("I8").Value = 0.02710117482667488
$ws.Range("J8").Value = 0.02710117482667489
$ws.Range("M8").Value = 10.34761366666667
$ws.Range("N8").Value = 31.042841
$ws.Range("O8").Value = 0.2299953477621856
$ws.Range("P8").Value = 0.2299953477621856
$ws.Range("Q8").Value = 7.232143796293
$ws.Range("R8").Value = 65.089294166637
$ws.Range("S8").Value = 0.006233144129024879
$ws.Range("T8").Value = 0.006233144129024881

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.6989190000000001
$ws.Range("H9").Value = 2.096757
$ws.Range("I9").Value = 0.02710117482667488
$ws.Range("J9").Value = 0.02710117482667489
$ws.Range("O9").Value = 0.6794731949692173
$ws.Range("P9").Value = 0.6794731949692174
$ws.Range("Q9").Value = 21.365857612151
$ws.Range("R9").Value = 192.292718509359
$ws.Range("S9").Value = 0.01841452184690011
$ws.Range("T9").Value = 0.01841452184690011

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.6989190000000001
$ws.Range("H10").Value = 2.096757
$ws.Range("I10").Value = 0.02710117482667488
$ws.Range("J10").Value = 0.02710117482667489
$ws.Range("M10").Value = 4.073058666666666
$ws.Range("N10").Value = 12.219176
$ws.Range("O10").Value = 0.09053145726859702
$ws.Range("P10").Value = 0.09053145726859703
$ws.Range("Q10").Value = 2.846738090248
$ws.Range("R10").Value = 25.620642812232
$ws.Range("S10").Value = 0.002453508850749894
$ws.Range("T10").Value = 0.002453508850749895

